# Dodanie podziału treningu na części
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F header ("Trening"), formatted like the other header cells
$ws.Range("F1").Value = "Trening"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Fill in the full data table (rows 2-13, columns A-F)
$ws.Cells.Item(2, 1).Value = 45685.64626446759
$ws.Cells.Item(2, 2).Value = 952.2
$ws.Cells.Item(2, 3).Value = 12.8
$ws.Cells.Item(2, 4).Value = 3.263509546007428
$ws.Cells.Item(2, 5).Value = "10-15"
$ws.Cells.Item(2, 6).Value = "Duża Gra"

$ws.Cells.Item(3, 1).Value = 45685.65125289352
$ws.Cells.Item(3, 2).Value = 1383.2
$ws.Cells.Item(3, 3).Value = 11.9
$ws.Cells.Item(3, 4).Value = 3.066205535616193
$ws.Cells.Item(3, 5).Value = "10-15"
$ws.Cells.Item(3, 6).Value = "Duża Gra"

$ws.Cells.Item(4, 1).Value = 45685.66516840277
$ws.Cells.Item(4, 2).Value = 2585.5
$ws.Cells.Item(4, 3).Value = 14.86
$ws.Cells.Item(4, 4).Value = 3.276531628199986
$ws.Cells.Item(4, 5).Value = "10-15"
$ws.Cells.Item(4, 6).Value = "Duża Gra"

$ws.Cells.Item(5, 1).Value = 45685.64626099537
$ws.Cells.Item(5, 2).Value = 951.9
$ws.Cells.Item(5, 3).Value = 9.49
$ws.Cells.Item(5, 4).Value = 3.115101371492657
$ws.Cells.Item(5, 5).Value = "5-10"
$ws.Cells.Item(5, 6).Value = "Duża Gra"

$ws.Cells.Item(6, 1).Value = 45685.6512505787
$ws.Cells.Item(6, 2).Value = 1383
$ws.Cells.Item(6, 3).Value = 9.96
$ws.Cells.Item(6, 4).Value = 2.997194494519916
$ws.Cells.Item(6, 5).Value = "5-10"
$ws.Cells.Item(6, 6).Value = "Duża Gra"

$ws.Cells.Item(7, 1).Value = 45685.65518344907
$ws.Cells.Item(7, 2).Value = 1722.8
$ws.Cells.Item(7, 3).Value = 9.26
$ws.Cells.Item(7, 4).Value = 2.782028470720563
$ws.Cells.Item(7, 5).Value = "5-10"
$ws.Cells.Item(7, 6).Value = "Duża Gra"

$ws.Cells.Item(8, 1).Value = 45685.66895081018
$ws.Cells.Item(8, 2).Value = 2912.3
$ws.Cells.Item(8, 3).Value = 12.08
$ws.Cells.Item(8, 4).Value = 3.491810968944004
$ws.Cells.Item(8, 5).Value = "10-15"
$ws.Cells.Item(8, 6).Value = "Mała Gra"

$ws.Cells.Item(9, 1).Value = 45685.6752667824
$ws.Cells.Item(9, 2).Value = 3458
$ws.Cells.Item(9, 3).Value = 14.66
$ws.Cells.Item(9, 4).Value = 3.470280102321081
$ws.Cells.Item(9, 5).Value = "10-15"
$ws.Cells.Item(9, 6).Value = "Mała Gra"

$ws.Cells.Item(10, 1).Value = 45685.67777256944
$ws.Cells.Item(10, 2).Value = 3674.5
$ws.Cells.Item(10, 3).Value = 11.96
$ws.Cells.Item(10, 4).Value = 3.452535833631243
$ws.Cells.Item(10, 5).Value = "10-15"
$ws.Cells.Item(10, 6).Value = "Mała Gra"

$ws.Cells.Item(11, 1).Value = 45685.66812326389
$ws.Cells.Item(11, 2).Value = 2840.8
$ws.Cells.Item(11, 3).Value = 9.67
$ws.Cells.Item(11, 4).Value = 2.932180132184711
$ws.Cells.Item(11, 5).Value = "5-10"
$ws.Cells.Item(11, 6).Value = "Mała Gra"

$ws.Cells.Item(12, 1).Value = 45685.66894849537
$ws.Cells.Item(12, 2).Value = 2912.1
$ws.Cells.Item(12, 3).Value = 9.1
$ws.Cells.Item(12, 4).Value = 3.239083426339286
$ws.Cells.Item(12, 5).Value = "5-10"
$ws.Cells.Item(12, 6).Value = "Mała Gra"

$ws.Cells.Item(13, 1).Value = 45685.67777025463
$ws.Cells.Item(13, 2).Value = 3674.3
$ws.Cells.Item(13, 3).Value = 8.61
$ws.Cells.Item(13, 4).Value = 2.997845990317207
$ws.Cells.Item(13, 5).Value = "5-10"
$ws.Cells.Item(13, 6).Value = "Mała Gra"

# Apply date/time display format to column A.
# First cell registers both the lower-case and upper-case format codes,
# the rest reuse the already-registered upper-case format.
$ws.Cells.Item(2, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(2, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(6, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(7, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(9, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(12, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "Applied training split edit"
